$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Russell Westbrook", "PG,SG", "Denver Nuggets"),
    @("Brandin Podziemski", "PG,SG", "Golden State Warriors"),
    @("Stephon Castle", "PG,SG", "San Antonio Spurs"),
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Pascal Siakam", "SF,PF,C", "Indiana Pacers"),
    @("Kyle Filipowski", "PF,C", "Utah Jazz"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Andre Drummond", "C", "Philadelphia 76ers"),
    @("Ty Jerome", "PG,SG", "Cleveland Cavaliers"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Jakob Poeltl", "C", "Toronto Raptors"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Payton Pritchard", "PG,SG", "Boston Celtics")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
